$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# ---------------------------------------------------------------------------
# 1. Bump the IG "Date" metadata value (Property/Value table on Metadata).
# ---------------------------------------------------------------------------
$wsMeta.Range("B8").Value = "2025-07-24T13:17:05+00:00"

# ---------------------------------------------------------------------------
# 2. Add a new row describing the "exerciceProfessionnel" element at the
#    bottom of the Elements table, mirroring row 10's layout/formatting.
# ---------------------------------------------------------------------------

# Copy row 10's formatting (fill/border/alignment) onto the new row 11 first.
$wsElem.Range("A10:AJ10").Copy()
$wsElem.Range("A11:AJ11").PasteSpecial(-4122)

# Copy row 10's values onto row 11 too, so every column (including the many
# "blank" ones) ends up populated/typed exactly like an existing data row.
$wsElem.Range("A10:AJ10").Copy()
$wsElem.Range("A11:AJ11").PasteSpecial(-4163)

# Min/Max (and Base Min/Base Max) for this element are both "1" - row 10 only
# has a ready-made text "1" in columns G/AH, so borrow those into F11/AG11
# (G11/AH11 already read "1" after the row copy above).
$wsElem.Range("G10").Copy()
$wsElem.Range("F11").PasteSpecial(-4163)
$wsElem.Range("AH10").Copy()
$wsElem.Range("AG11").PasteSpecial(-4163)

# Now overwrite the cells that are specific to this new element.
$wsElem.Range("A11").Value = "SituationOperationnelle.exerciceProfessionnel"
$wsElem.Range("B11").Value = "SituationOperationnelle.exerciceProfessionnel"
$wsElem.Range("K11").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/ExerciceProfessionnel)`n"
$wsElem.Range("L11").Value = "Lien vers la classe ExerciceProfessionnel."
$wsElem.Range("M11").Value = "Lien vers la classe ExerciceProfessionnel."
$wsElem.Range("AF11").Value = "SituationOperationnelle.exerciceProfessionnel"

# Entering a value that ends with a newline makes the grid auto-grow the row
# height; AutoFit puts it back to the sheet's normal (default) row height.
$wsElem.Rows.Item(11).AutoFit()

# Widen column K (Type(s)) to fit the new, longer reference text.
$wsElem.Columns.Item(11).ColumnWidth = 74.6
